# Update the cryptos price list: refresh Price (D) and Volume(1h) (E) columns.
# Values that look like plain decimal numbers (e.g. "214.76") are forced to
# text via NumberFormat "@" so they stay as strings instead of being
# auto-converted to numeric values by Excel, matching the source data which
# stores these as text (prices use "." as a thousands separator in some rows,
# e.g. "28.403.31", so they must remain text throughout).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.403.31"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").Value = "1.595.84"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.76"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.495"
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.09"
$ws.Range("E8").Value = "  +8.74%  "
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "1.822.09"
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("D13").Value = "1.595.12"
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.533"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "28.400.44"
$ws.Range("E16").Value = "  +4.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.20"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.45"
$ws.Range("E18").Value = "  +4.78%  "
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.51"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.33"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.88"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.22"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.108"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "1.398.50"
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("E36").Value = "  -5.15%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.52"
$ws.Range("E39").Value = "  +7.74%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.816"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.73"
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.89"
$ws.Range("E44").Value = "  +7.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.985"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.53"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "1.731.97"
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.61"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0526"
$ws.Range("E51").Value = "  +0.15%  "
